# Apply the "effort" log update:
#  - Row 16: merge Additional Effort [h] (C16=0.5) into Effort [h] (B16),
#            giving B16 = 5.5, and clear C16.
#  - Append a new log entry as row 32:
#       Date = 2012-10-26 (serial 41208), Effort [h] = 4,
#       Comment = "Manual continued, widely completed as draft but
#                  without chapter Use Cases"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: fold the additional-effort hours into the main effort cell ---
$ws.Range("B16").Value = 5.5
$ws.Range("C16").Value = ""

# --- New row 32: next day of work on the manual ---
$ws.Range("A32").Value = 41208
$ws.Range("B32").Value = 4
$ws.Range("D32").Value = "Manual continued, widely completed as draft but without chapter Use Cases"
